# Pumped Hydro Example - update turbine efficiency and add an
# "Olympic swimming pool" comparison for the 5 MWh storage scenario.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update turbine efficiency assumption (B4): 0.72 -> 0.9 -------------
# This ripples through B12 (mass -> volume conversion) automatically.
$ws.Range("B4").Value = 0.9

# --- New rows comparing 5 MWh of pumped storage to Olympic pools --------
# Row 13: reference constant - volume of an Olympic swimming pool (m^3)
$ws.Range("C13").Value = "m^3"
$ws.Range("A15").Value = "5 MWh storage equals:"
$ws.Range("A16").Value = "Olympic swimming pools"
$ws.Range("A13").Value = "volume of olympic swimming pool"

$ws.Range("B13").Formula = "=2500"

# Row 15: storage volume required for 5 MWh, scaled off the B12 result
$ws.Range("B15").Formula = "=B12/B6*5000"
$ws.Range("C15").Value = "m^3"

# Row 16: that volume expressed in Olympic swimming pools
$ws.Range("B16").Formula = "=SUM(B15)"

# Leave selection on the newly added total cell, matching the edited file.
$ws.Range("B16").Select()
